# Auto-generated Excel COM-interop script
# Applies the numeric cell-value updates described by the target diff
# for Sheets/Coeurl_Profits.xlsx (tabs: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1173.8889
$ws.Range("I31").Value = 1173.8889
$ws.Range("K31").Value = 3521.6667
$ws.Range("M31").Value = -3291.6667

$ws.Range("H125").Value = 3558.0588
$ws.Range("J125").Value = 3655.5715
$ws.Range("L125").Value = 32900.1435
$ws.Range("N125").Value = -37820.1435

$ws.Range("H132").Value = 1167.1132
$ws.Range("I132").Value = 1134.4509
$ws.Range("K132").Value = 3403.3527
$ws.Range("M132").Value = -873.3527000000004

$ws.Range("H137").Value = 1745.6765
$ws.Range("I137").Value = 1583.7778
$ws.Range("K137").Value = 4751.3334
$ws.Range("M137").Value = -2201.3334

$ws.Range("H138").Value = 6852202
$ws.Range("I138").Value = 1327.7142
$ws.Range("K138").Value = 3983.1426
$ws.Range("M138").Value = 1156.8574

$ws.Range("H141").Value = 862.8889
$ws.Range("I141").Value = 862.8889
$ws.Range("K141").Value = 2588.6667
$ws.Range("M141").Value = 2591.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8407.053
$ws.Range("I32").Value = 3757.9827
$ws.Range("K32").Value = 3757.9827
$ws.Range("M32").Value = -3470.9827

$ws.Range("H61").Value = 4096.8037
$ws.Range("I61").Value = 2916.8064
$ws.Range("K61").Value = 2916.8064
$ws.Range("M61").Value = -2704.8064

$ws.Range("H74").Value = 8870.885
$ws.Range("I74").Value = 1654.7142
$ws.Range("K74").Value = 1654.7142
$ws.Range("M74").Value = -780.7141999999999

$ws.Range("H77").Value = 8870.885
$ws.Range("I77").Value = 1654.7142
$ws.Range("K77").Value = 8273.571
$ws.Range("M77").Value = -3905.571

$ws.Range("H92").Value = 60000
$ws.Range("J92").Value = 60000
$ws.Range("L92").Value = 60000
$ws.Range("N92").Value = -64992

$ws.Range("H96").Value = 30562.334
$ws.Range("J96").Value = 30562.334
$ws.Range("L96").Value = 30562.334
$ws.Range("N96").Value = -36054.334

$ws.Range("H110").Value = 4266.0884
$ws.Range("I110").Value = 4581.6333
$ws.Range("J110").Value = 1899.5
$ws.Range("K110").Value = 4581.6333
$ws.Range("L110").Value = 1899.5
$ws.Range("M110").Value = -2536.6333
$ws.Range("N110").Value = -5989.5

$ws.Range("H122").Value = 1329
$ws.Range("I122").Value = 1172.2941
$ws.Range("K122").Value = 3516.8823
$ws.Range("M122").Value = -1066.8823

$ws.Range("H136").Value = 4096.8037
$ws.Range("I136").Value = 2916.8064
$ws.Range("K136").Value = 8750.4192
$ws.Range("M136").Value = -6200.4192

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2091.7693
$ws.Range("I105").Value = 2219.9
$ws.Range("J105").Value = 1664.6666
$ws.Range("K105").Value = 2219.9
$ws.Range("L105").Value = 1664.6666
$ws.Range("M105").Value = -472.9000000000001
$ws.Range("N105").Value = -5158.6666

$ws.Range("H132").Value = 25000
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 1888.0941
$ws.Range("I134").Value = 1882.4073
$ws.Range("J134").Value = 2003.25
$ws.Range("K134").Value = 5647.2219
$ws.Range("L134").Value = 6009.75
$ws.Range("M134").Value = -3112.2219
$ws.Range("N134").Value = -11079.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 50000
$ws.Range("I56").Value = 50000
$ws.Range("K56").Value = 50000
$ws.Range("M56").Value = -49155

$ws.Range("H58").Value = 3145.3264
$ws.Range("I58").Value = 2791.366
$ws.Range("K58").Value = 2791.366
$ws.Range("M58").Value = -2588.366

$ws.Range("H93").Value = 9678.25
$ws.Range("I93").Value = 5489.4287
$ws.Range("K93").Value = 5489.4287
$ws.Range("M93").Value = -3617.4287

$ws.Range("H103").Value = 28331.666
$ws.Range("I103").Value = 25998
$ws.Range("K103").Value = 25998
$ws.Range("M103").Value = -24826

$ws.Range("H132").Value = 2899.375
$ws.Range("J132").Value = 3448.1428
$ws.Range("L132").Value = 10344.4284
$ws.Range("N132").Value = -15404.4284

$ws.Range("H136").Value = 3145.3264
$ws.Range("I136").Value = 2791.366
$ws.Range("K136").Value = 8374.098
$ws.Range("M136").Value = -5824.098

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 577.9091
$ws.Range("I5").Value = 535.75
$ws.Range("K5").Value = 1607.25
$ws.Range("M5").Value = -1495.25

$ws.Range("H39").Value = 7289.622
$ws.Range("J39").Value = 9227.272000000001
$ws.Range("L39").Value = 27681.816
$ws.Range("N39").Value = -28269.816

$ws.Range("H55").Value = 2046.1538
$ws.Range("J55").Value = 3214.1428
$ws.Range("L55").Value = 9642.428400000001
$ws.Range("N55").Value = -9996.428400000001

$ws.Range("H110").Value = 20013.5
$ws.Range("I110").Value = 20013.5
$ws.Range("K110").Value = 60040.5
$ws.Range("M110").Value = -55950.5

$ws.Range("H135").Value = 577.9091
$ws.Range("I135").Value = 535.75
$ws.Range("K135").Value = 4821.75
$ws.Range("M135").Value = -2286.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3748
$ws.Range("I80").Value = 3494
$ws.Range("K80").Value = 3494
$ws.Range("M80").Value = -2496

$ws.Range("H83").Value = 3748
$ws.Range("I83").Value = 3494
$ws.Range("K83").Value = 17470
$ws.Range("M83").Value = -12478

$ws.Range("H97").Value = 1040.25
$ws.Range("J97").Value = 2171.25
$ws.Range("L97").Value = 2171.25
$ws.Range("N97").Value = -3163.25

$ws.Range("H113").Value = 2451.3333
$ws.Range("J113").Value = 1430.6666
$ws.Range("L113").Value = 1430.6666
$ws.Range("N113").Value = -5770.6666

$ws.Range("H122").Value = 2452.2
$ws.Range("I122").Value = 1795.4762
$ws.Range("K122").Value = 5386.4286
$ws.Range("M122").Value = -2936.4286

$ws.Range("H132").Value = 2586.8936
$ws.Range("I132").Value = 2274.8647
$ws.Range("J132").Value = 3741.4
$ws.Range("K132").Value = 6824.5941
$ws.Range("L132").Value = 11224.2
$ws.Range("M132").Value = -4294.5941
$ws.Range("N132").Value = -16284.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4670.3687
$ws.Range("I40").Value = 3513.0715
$ws.Range("J40").Value = 7910.8
$ws.Range("K40").Value = 3513.0715
$ws.Range("L40").Value = 7910.8
$ws.Range("M40").Value = -3377.0715
$ws.Range("N40").Value = -8182.8

$ws.Range("H60").Value = 80000
$ws.Range("J60").Value = 80000
$ws.Range("L60").Value = 80000
$ws.Range("N60").Value = -81018

$ws.Range("H61").Value = 4097.8125
$ws.Range("I61").Value = 4147.5
$ws.Range("J61").Value = 3750
$ws.Range("K61").Value = 4147.5
$ws.Range("L61").Value = 3750
$ws.Range("M61").Value = -3945.5
$ws.Range("N61").Value = -4154

$ws.Range("H113").Value = 4097.8125
$ws.Range("I113").Value = 4147.5
$ws.Range("J113").Value = 3750
$ws.Range("K113").Value = 4147.5
$ws.Range("L113").Value = 3750
$ws.Range("M113").Value = -1977.5
$ws.Range("N113").Value = -8090

$ws.Range("H132").Value = 4417.2856
$ws.Range("I132").Value = 4272
$ws.Range("K132").Value = 12816
$ws.Range("M132").Value = -10286

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 38333
$ws.Range("J33").Value = 38333
$ws.Range("N33").Value = -38833

$ws.Range("H36").Value = 38333
$ws.Range("J36").Value = 38333
$ws.Range("N36").Value = -38833

$ws.Range("H37").Value = 22295
$ws.Range("J37").Value = 22295
$ws.Range("L37").Value = 22295
$ws.Range("N37").Value = -22701

$ws.Range("H94").Value = 11875.556
$ws.Range("I94").Value = 18499.5
$ws.Range("J94").Value = 9983
$ws.Range("K94").Value = 18499.5
$ws.Range("L94").Value = 9983
$ws.Range("M94").Value = -17598.5
$ws.Range("N94").Value = -11785

$ws.Range("H113").Value = 1207.5862
$ws.Range("I113").Value = 1110.4445
$ws.Range("K113").Value = 3331.3335
$ws.Range("M113").Value = -1161.3335

$ws.Range("H132").Value = 1507.3478
$ws.Range("I132").Value = 1346.0328
$ws.Range("K132").Value = 4038.0984
$ws.Range("M132").Value = -1508.0984

$ws.Range("H136").Value = 1986.16
$ws.Range("I136").Value = 1757.0454
$ws.Range("K136").Value = 5271.1362
$ws.Range("M136").Value = -2721.1362
